$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new withdrawal records ("prelievi") were added to the log on
# 26/04/2018 and 27/04/2018 (rows 10 and 11), following the same layout
# as every other row: Data | Prelevante | Materiale | UnitaMisura | Quantità

# Replicate the formatting of the last existing data row (row 9) onto the
# two new rows so the new cells render with the same date / text / number
# styles as the rest of the table.
for ($col = 1; $col -le 5; $col++) {
    $ws.Cells.Item(9, $col).Copy($ws.Cells.Item(10, $col))
    $ws.Cells.Item(9, $col).Copy($ws.Cells.Item(11, $col))
}

# Row 10: 26/04/2018 - Saldakeeva Elena - Adesivo Leggero Bianco - Mt. - 1
$ws.Cells.Item(10, 1).Value = 43216
$ws.Cells.Item(10, 2).Value = "Saldakeeva Elena"
$ws.Cells.Item(10, 3).Value = "Adesivo Leggero Bianco"
$ws.Cells.Item(10, 4).Value = "Mt."
$ws.Cells.Item(10, 5).Value = 1

# Row 11: 27/04/2018 - Katia D'Alesio - Shopping Bag - N°. - 40
$ws.Cells.Item(11, 1).Value = 43217
$ws.Cells.Item(11, 2).Value = "Katia D'Alesio"
$ws.Cells.Item(11, 3).Value = "Shopping Bag"
$ws.Cells.Item(11, 4).Value = "N°."
$ws.Cells.Item(11, 5).Value = 40
